# Nexial step-showcase.xlsx update:
#  - insert a new "aws.ses" command group into the '#system' sheet, which acts
#    as a lookup table for data-validation driven command pickers.
#  - "target" (column A) lists all group names; "aws.ses" needs to be spliced
#    in alphabetically between "aws.s3" and "base" (row 3), pushing every
#    following target row down by one.
#  - a new column C holds the "aws.ses" command list (header + 2 commands);
#    every existing group column from the old C ("base") through old Z ("xml")
#    shifts one column to the right (D..AA).
#  - the defined names that point at those columns must be updated to match,
#    and a new "aws.ses" defined name is added pointing at the new column C.

function ColLetter([int]$n) {
    $letter = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $letter = [char](65 + $rem) + $letter
        $n = [int](($n - $rem - 1) / 26)
    }
    return $letter
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. Read the entire existing table (A1:Z117) in one shot.
# ---------------------------------------------------------------------------
$oldLastRow = 117
$oldLastCol = 26   # Z
$srcAddr = "A1:" + (ColLetter $oldLastCol) + $oldLastRow
$srcRange = $ws.Range($srcAddr)
$src = $srcRange.Value2

$newLastRow = $oldLastRow
$newLastCol = $oldLastCol + 1   # AA (one new column inserted)

# ---------------------------------------------------------------------------
# 2. Build the new grid in memory.
# ---------------------------------------------------------------------------
$dst = New-Object 'object[,]' $newLastRow, $newLastCol

# -- 2a. Column A ("target"): copy rows 1-2 as-is, then splice "aws.ses" in
#        at row 3, shifting old rows 3..26 down to 4..27.
$dst[0, 0] = $src[1, 1]
$dst[1, 0] = $src[2, 1]
$dst[2, 0] = "aws.ses"
for ($r = 3; $r -le $oldLastRow; $r++) {
    $v = $src[$r, 1]
    if ($v -ne $null) {
        $dst[$r, 0] = $v
    }
}

# -- 2b. Column B ("aws.s3"): unchanged, straight copy.
for ($r = 1; $r -le $oldLastRow; $r++) {
    $v = $src[$r, 2]
    if ($v -ne $null) {
        $dst[$r - 1, 1] = $v
    }
}

# -- 2c. Old columns C..Z (3..26) shift right by one -> new columns D..AA (4..27).
for ($c = 3; $c -le $oldLastCol; $c++) {
    for ($r = 1; $r -le $oldLastRow; $r++) {
        $v = $src[$r, $c]
        if ($v -ne $null) {
            $dst[$r - 1, $c] = $v
        }
    }
}

# -- 2d. New column C holds the "aws.ses" command group.
$dst[0, 2] = "aws.ses"
$dst[1, 2] = "sendMail(profile,to,subject,body)"
$dst[2, 2] = "sendTextMail(profile,to,subject,body)"

# ---------------------------------------------------------------------------
# 3. Write the new grid back in one shot, then write it out to the sheet.
# ---------------------------------------------------------------------------
$dstAddr = "A1:" + (ColLetter $newLastCol) + $newLastRow
$destRange = $ws.Range($dstAddr)
$destRange.Value2 = $dst

# ---------------------------------------------------------------------------
# 4. Update defined names that referred to the shifted columns, and add the
#    new "aws.ses" defined name.
# ---------------------------------------------------------------------------
$shiftedNames = @{
    "base"      = "`$D`$2:`$D`$36";
    "csv"       = "`$E`$2:`$E`$5";
    "desktop"   = "`$F`$2:`$F`$92";
    "excel"     = "`$G`$2:`$G`$14";
    "external"  = "`$H`$2:`$H`$3";
    "image"     = "`$I`$2:`$I`$5";
    "io"        = "`$J`$2:`$J`$24";
    "jms"       = "`$K`$2:`$K`$4";
    "json"      = "`$L`$2:`$L`$14";
    "mail"      = "`$M`$2:`$M`$2";
    "number"    = "`$N`$2:`$N`$15";
    "pdf"       = "`$O`$2:`$O`$16";
    "rdbms"     = "`$P`$2:`$P`$7";
    "redis"     = "`$Q`$2:`$Q`$10";
    "sms"       = "`$R`$2:`$R`$2";
    "sound"     = "`$S`$2:`$S`$5";
    "ssh"       = "`$T`$2:`$T`$9";
    "step"      = "`$U`$2:`$U`$4";
    "target"    = "`$A`$2:`$A`$27";
    "web"       = "`$V`$2:`$V`$117";
    "webalert"  = "`$W`$2:`$W`$8";
    "webcookie" = "`$X`$2:`$X`$8";
    "ws"        = "`$Y`$2:`$Y`$17";
    "ws.async"  = "`$Z`$2:`$Z`$8";
    "xml"       = "`$AA`$2:`$AA`$11";
}

foreach ($name in $shiftedNames.Keys) {
    $ref = "='#system'!" + $shiftedNames[$name]
    $wb.Names.Item($name).RefersTo = $ref
}

$wb.Names.Add("aws.ses", "='#system'!`$C`$2:`$C`$3")

Write-Host "step-showcase.xlsx updated: aws.ses command group inserted"
